$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(39, 1).Value = 37
$ws.Cells.Item(39, 2).Value = 2.686192441575164
$ws.Cells.Item(39, 3).Value = 4940.8
$ws.Cells.Item(39, 4).Value = 0.01922476690873743
$ws.Cells.Item(39, 5).Value = 33
$ws.Cells.Item(39, 6).Value = 215.4
$ws.Cells.Item(39, 7).Value = "Portgas D Åce "
$ws.Cells.Item(39, 8).Value = "SOLO"
$ws.Cells.Item(39, 9).Value = 0.1187514263010913
$ws.Cells.Item(39, 10).Value = 6.2
$ws.Cells.Item(39, 11).Value = 0.003086678567510356

$ws.Cells.Item(40, 1).Value = 38
$ws.Cells.Item(40, 2).Value = 6.764476155827955
$ws.Cells.Item(40, 3).Value = 13106.2
$ws.Cells.Item(40, 4).Value = 0.01670651213554389
$ws.Cells.Item(40, 5).Value = 32.8
$ws.Cells.Item(40, 6).Value = 144.2
$ws.Cells.Item(40, 7).Value = "BigFather Rengar"
$ws.Cells.Item(40, 8).Value = "SOLO"
$ws.Cells.Item(40, 9).Value = 0.07387201437406728
$ws.Cells.Item(40, 10).Value = 10.6
$ws.Cells.Item(40, 11).Value = 0.00546712006355476

$ws.Cells.Item(41, 1).Value = 39
$ws.Cells.Item(41, 2).Value = 4.601146283614828
$ws.Cells.Item(41, 3).Value = 8222
$ws.Cells.Item(41, 4).Value = 0.03186341646984259
$ws.Cells.Item(41, 5).Value = 55.6
$ws.Cells.Item(41, 6).Value = 280.8
$ws.Cells.Item(41, 7).Value = "Cevahir Akkanat"
$ws.Cells.Item(41, 8).Value = "SOLO"
$ws.Cells.Item(41, 9).Value = 0.154526359571592
$ws.Cells.Item(41, 10).Value = 13
$ws.Cells.Item(41, 11).Value = 0.006825017206633095

$ws.Cells.Item(42, 1).Value = 40
$ws.Cells.Item(42, 2).Value = 2.66602457655264
$ws.Cells.Item(42, 3).Value = 8027.4
$ws.Cells.Item(42, 4).Value = 0.01693789438724676
$ws.Cells.Item(42, 5).Value = 51
$ws.Cells.Item(42, 6).Value = 193.6
$ws.Cells.Item(42, 7).Value = "Jimmy L"
$ws.Cells.Item(42, 8).Value = "SOLO"
$ws.Cells.Item(42, 9).Value = 0.06429757555629359
$ws.Cells.Item(42, 10).Value = 5
$ws.Cells.Item(42, 11).Value = 0.001660577881102624

$ws.Cells.Item(43, 1).Value = 41
$ws.Cells.Item(43, 2).Value = 4.709765260580623
$ws.Cells.Item(43, 3).Value = 8880
$ws.Cells.Item(43, 4).Value = 0.03577748863238895
$ws.Cells.Item(43, 5).Value = 67.6
$ws.Cells.Item(43, 6).Value = 204.2
$ws.Cells.Item(43, 7).Value = "QUDURMAZSAN400RP"
$ws.Cells.Item(43, 8).Value = "SOLO"
$ws.Cells.Item(43, 9).Value = 0.1078900314795383
$ws.Cells.Item(43, 10).Value = 8.4
$ws.Cells.Item(43, 11).Value = 0.004447037425673313

$ws.Cells.Item(44, 1).Value = 42
$ws.Cells.Item(44, 2).Value = 6.322405523740514
$ws.Cells.Item(44, 3).Value = 11499.8
$ws.Cells.Item(44, 4).Value = 0.05329568716167528
$ws.Cells.Item(44, 5).Value = 97.2
$ws.Cells.Item(44, 6).Value = 286.4
$ws.Cells.Item(44, 7).Value = "OsmanGazi0505"
$ws.Cells.Item(44, 8).Value = "DUO_CARRY"
$ws.Cells.Item(44, 9).Value = 0.1618235216992612
$ws.Cells.Item(44, 10).Value = 16.8
$ws.Cells.Item(44, 11).Value = 0.009080032401669283

$ws.Cells.Item(45, 1).Value = 43
$ws.Cells.Item(45, 2).Value = 2.04804347826087
$ws.Cells.Item(45, 3).Value = 3512.2
$ws.Cells.Item(45, 4).Value = 0.03210989430096982
$ws.Cells.Item(45, 5).Value = 52.6
$ws.Cells.Item(45, 6).Value = 147
$ws.Cells.Item(45, 7).Value = "Mallorn"
$ws.Cells.Item(45, 8).Value = "DUO_CARRY"
$ws.Cells.Item(45, 9).Value = 0.08620600414078675
$ws.Cells.Item(45, 10).Value = 12.2
$ws.Cells.Item(45, 11).Value = 0.006946169772256728

$ws.Cells.Item(46, 1).Value = 44
$ws.Cells.Item(46, 2).Value = 2.122644163150492
$ws.Cells.Item(46, 3).Value = 3018.4
$ws.Cells.Item(46, 4).Value = 0.01828410689170183
$ws.Cells.Item(46, 5).Value = 26
$ws.Cells.Item(46, 6).Value = 95.6
$ws.Cells.Item(46, 7).Value = "TheImperium"
$ws.Cells.Item(46, 8).Value = "SOLO"
$ws.Cells.Item(46, 9).Value = 0.06722925457102671
$ws.Cells.Item(46, 10).Value = 2
$ws.Cells.Item(46, 11).Value = 0.001406469760900141

$ws.Cells.Item(47, 1).Value = 45
$ws.Cells.Item(47, 2).Value = 5.449437357377095
$ws.Cells.Item(47, 3).Value = 10147.6
$ws.Cells.Item(47, 4).Value = 0.03793728796143418
$ws.Cells.Item(47, 5).Value = 72
$ws.Cells.Item(47, 6).Value = 323.6
$ws.Cells.Item(47, 7).Value = "Memoşata"
$ws.Cells.Item(47, 8).Value = "SOLO"
$ws.Cells.Item(47, 9).Value = 0.17855598675708
$ws.Cells.Item(47, 10).Value = 11.6
$ws.Cells.Item(47, 11).Value = 0.006341493230895424

$ws.Cells.Item(48, 1).Value = 46
$ws.Cells.Item(48, 2).Value = 2.720225988700565
$ws.Cells.Item(48, 3).Value = 4814.8
$ws.Cells.Item(48, 4).Value = 0.01751412429378531
$ws.Cells.Item(48, 5).Value = 31
$ws.Cells.Item(48, 6).Value = 136.6
$ws.Cells.Item(48, 7).Value = "god damn u suck"
$ws.Cells.Item(48, 8).Value = "SOLO"
$ws.Cells.Item(48, 9).Value = 0.07717514124293785
$ws.Cells.Item(48, 10).Value = 4.4
$ws.Cells.Item(48, 11).Value = 0.002485875706214689

$ws.Cells.Item(49, 1).Value = 47
$ws.Cells.Item(49, 2).Value = 1.156962025316456
$ws.Cells.Item(49, 3).Value = 1645.2
$ws.Cells.Item(49, 4).Value = 0.02165963431786217
$ws.Cells.Item(49, 5).Value = 30.8
$ws.Cells.Item(49, 6).Value = 249
$ws.Cells.Item(49, 7).Value = "hobosapien1"
$ws.Cells.Item(49, 8).Value = "DUO_SUPPORT"
$ws.Cells.Item(49, 9).Value = 0.1751054852320675
$ws.Cells.Item(49, 10).Value = 3.6
$ws.Cells.Item(49, 11).Value = 0.002531645569620253

$ws.Cells.Item(50, 1).Value = 48
$ws.Cells.Item(50, 2).Value = 4.577926580861201
$ws.Cells.Item(50, 3).Value = 8986.6
$ws.Cells.Item(50, 4).Value = 0.01983526296622257
$ws.Cells.Item(50, 5).Value = 37
$ws.Cells.Item(50, 6).Value = 366.6
$ws.Cells.Item(50, 7).Value = "PYKEntakill"
$ws.Cells.Item(50, 8).Value = "DUO_SUPPORT"
$ws.Cells.Item(50, 9).Value = 0.1949690809435287
$ws.Cells.Item(50, 10).Value = 16.4
$ws.Cells.Item(50, 11).Value = 0.008763317485602414

$ws.Cells.Item(51, 1).Value = 49
$ws.Cells.Item(51, 2).Value = 2.732152492161565
$ws.Cells.Item(51, 3).Value = 4754.8
$ws.Cells.Item(51, 4).Value = 0.01058489950476475
$ws.Cells.Item(51, 5).Value = 18.2
$ws.Cells.Item(51, 6).Value = 286.2
$ws.Cells.Item(51, 7).Value = "Portgas D Åce "
$ws.Cells.Item(51, 8).Value = "SOLO"
$ws.Cells.Item(51, 9).Value = 0.1699525553653752
$ws.Cells.Item(51, 10).Value = 21.2
$ws.Cells.Item(51, 11).Value = 0.01215288164343084

$ws.Cells.Item(52, 1).Value = 50
$ws.Cells.Item(52, 2).Value = 3.151435406698565
$ws.Cells.Item(52, 3).Value = 5269.2
$ws.Cells.Item(52, 4).Value = 0.02009569377990431
$ws.Cells.Item(52, 5).Value = 33.6
$ws.Cells.Item(52, 6).Value = 116
$ws.Cells.Item(52, 7).Value = "LS DUFFY"
$ws.Cells.Item(52, 8).Value = "SOLO"
$ws.Cells.Item(52, 9).Value = 0.06937799043062201
$ws.Cells.Item(52, 10).Value = 2.6
$ws.Cells.Item(52, 11).Value = 0.001555023923444976

$ws.Cells.Item(53, 1).Value = 51
$ws.Cells.Item(53, 2).Value = 6.383743300662566
$ws.Cells.Item(53, 3).Value = 12258.8
$ws.Cells.Item(53, 4).Value = 0.04023672888935409
$ws.Cells.Item(53, 5).Value = 81
$ws.Cells.Item(53, 6).Value = 253.6
$ws.Cells.Item(53, 7).Value = "BigFather Rengar"
$ws.Cells.Item(53, 8).Value = "SOLO"
$ws.Cells.Item(53, 9).Value = 0.13711238076999
$ws.Cells.Item(53, 10).Value = 14.8
$ws.Cells.Item(53, 11).Value = 0.007707535390238893

$ws.Cells.Item(54, 1).Value = 57
$ws.Cells.Item(54, 2).Value = 8.743415893791914
$ws.Cells.Item(54, 3).Value = 18485.6
$ws.Cells.Item(54, 4).Value = 0.04002778947189373
$ws.Cells.Item(54, 5).Value = 84
$ws.Cells.Item(54, 6).Value = 229.6
$ws.Cells.Item(54, 7).Value = "Châllénger "
$ws.Cells.Item(54, 8).Value = "SOLO"
$ws.Cells.Item(54, 9).Value = 0.1242981404229443
$ws.Cells.Item(54, 10).Value = 10
$ws.Cells.Item(54, 11).Value = 0.00490817992779415

$ws.Cells.Item(55, 1).Value = 58
$ws.Cells.Item(55, 2).Value = 2.732152492161565
$ws.Cells.Item(55, 3).Value = 4754.8
$ws.Cells.Item(55, 4).Value = 0.01058489950476475
$ws.Cells.Item(55, 5).Value = 18.2
$ws.Cells.Item(55, 6).Value = 286.2
$ws.Cells.Item(55, 7).Value = "Portgas D Åce "
$ws.Cells.Item(55, 8).Value = "SOLO"
$ws.Cells.Item(55, 9).Value = 0.1699525553653752
$ws.Cells.Item(55, 10).Value = 21.2
$ws.Cells.Item(55, 11).Value = 0.01215288164343084

$ws.Cells.Item(56, 1).Value = 59
$ws.Cells.Item(56, 2).Value = 3.151435406698565
$ws.Cells.Item(56, 3).Value = 5269.2
$ws.Cells.Item(56, 4).Value = 0.02009569377990431
$ws.Cells.Item(56, 5).Value = 33.6
$ws.Cells.Item(56, 6).Value = 116
$ws.Cells.Item(56, 7).Value = "LS DUFFY"
$ws.Cells.Item(56, 8).Value = "SOLO"
$ws.Cells.Item(56, 9).Value = 0.06937799043062201
$ws.Cells.Item(56, 10).Value = 2.6
$ws.Cells.Item(56, 11).Value = 0.001555023923444976

$ws.Cells.Item(57, 1).Value = 60
$ws.Cells.Item(57, 2).Value = 6.383743300662566
$ws.Cells.Item(57, 3).Value = 12258.8
$ws.Cells.Item(57, 4).Value = 0.04023672888935409
$ws.Cells.Item(57, 5).Value = 81
$ws.Cells.Item(57, 6).Value = 253.6
$ws.Cells.Item(57, 7).Value = "BigFather Rengar"
$ws.Cells.Item(57, 8).Value = "SOLO"
$ws.Cells.Item(57, 9).Value = 0.13711238076999
$ws.Cells.Item(57, 10).Value = 14.8
$ws.Cells.Item(57, 11).Value = 0.007707535390238893

$ws.Cells.Item(58, 1).Value = 61
$ws.Cells.Item(58, 2).Value = 1.901710291787398
$ws.Cells.Item(58, 3).Value = 4185.8
$ws.Cells.Item(58, 4).Value = 0.01349240180076153
$ws.Cells.Item(58, 5).Value = 29.6
$ws.Cells.Item(58, 6).Value = 204.4
$ws.Cells.Item(58, 7).Value = "MyDogaN"
$ws.Cells.Item(58, 8).Value = "DUO_SUPPORT"
$ws.Cells.Item(58, 9).Value = 0.0912805515236867
$ws.Cells.Item(58, 10).Value = 18.6
$ws.Cells.Item(58, 11).Value = 0.00816721733901452

# Apply the header-row style (bold, border, centered) to the new A-column cells
$ws.Range("A39").Copy()
$ws.Range("A45:A58").PasteSpecial(-4122)

